$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, border, alignment) from the
# existing last header cell (AC1) into the three new header cells so
# they match the style used by the rest of row 1.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

# Set the new header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every data row.
$ws.Range("AD2:AD51").Value = 75
$ws.Range("AE2:AE51").Value = 87
$ws.Range("AF2:AF51").Value = 0
